$d = $word.ActiveDocument

# --- Change 1: extend the "Layout" description paragraph (paragraph 2) with
# additional sentences describing the menu strip / text boxes. ---
$p2 = $d.Paragraphs.Item(2)
$p2.Range.InsertAfter(" The ")

$p2 = $d.Paragraphs.Item(2)
$p2.Range.InsertAfter("buttons and options are accessed through a menu strip on the top of the form")

$p2 = $d.Paragraphs.Item(2)
$p2.Range.InsertAfter(". From there, depending on the option, text boxes are given to enter the")

$p2 = $d.Paragraphs.Item(2)
$p2.Range.InsertAfter(" necessary information to add to the system.")

# --- Change 2: after the "-Defining the requirements..." paragraph (which is
# paragraph 5), add two new bullet-style paragraphs describing more steps,
# followed by one blank paragraph, before the "Screenshots:" heading. ---
$pDefining = $d.Paragraphs.Item(5)
$pDefining.Range.InsertParagraphAfter()

$pCreating = $d.Paragraphs.Item(6)
$pCreating.Range.InsertAfter("-Creating base layout of the menu strip and button dropdowns")

$pCreating = $d.Paragraphs.Item(6)
$pCreating.Range.InsertParagraphAfter()

$pGiveCode = $d.Paragraphs.Item(7)
$pGiveCode.Range.InsertAfter("-Give code to each drop down button in the menustrip")

$pGiveCode = $d.Paragraphs.Item(7)
$pGiveCode.Range.InsertParagraphAfter()
